# Update "Förändrad" date column (C2:C9) from 2023-10-08 (45207) to 2023-10-09 (45208)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

for ($row = 2; $row -le 9; $row++) {
    $ws.Cells.Item($row, 3).Value = 45208
}
